$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 13 - this shifts the existing "Step1" block
# (rows 13-22) down by one row, to rows 14-23.
$ws.Rows(13).Insert()

# New row 13: the "Step2" step that feeds the new (String[]) null value.
$ws.Range("D13").Value = "Step2"
$ws.Range("E13").Value = "'= (String[]) null"

# Append ".length" to the two existing step formulas (now on rows 15/16
# after the insert above).
$ws.Range("E15").Value = "'=flatten(null).length"
$ws.Range("E16").Value = "'=flatten(`$Step1).length"

# Row 17 was an empty spacer row before the edit; it now gets the new
# "result2" step.
$ws.Range("D17").Value = "result2"
$ws.Range("E17").Value = "'= flatten(`$Step2).length"

# Add the new "result2" column to the test table header (rows 21-22,
# shifted down from 20-21).
$ws.Range("G21").Value = "_res_.`$result2"
$ws.Range("G22").Value = "_res_.`$result2"

# Add the new test data row values (row 23, shifted down from 22).
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0

$ws.Range("M14").Select()
